# KPI_ACTION_TRACKER.xlsx - "Commit Dec 27th v2"
# Updates KPI statuses, refreshes the Appzen action-item wording with
# parenthesised dates, and scrolls/reselects the KPI sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# --- Status updates (column G = Progress) ---------------------------------
$ws.Range("G28").Value = "Complete"
$ws.Range("G29").Value = "Complete"
$ws.Range("G30").Value = "Overdue"

# --- Reword the Appzen action item (row 35, column F = How) ---------------
$ws.Range("F35").Value = "1. Appzen current configuration setting revisit (Jan 7, 2022)`n2. Public transporation configuration as a reference point (Jan 7, 2022)`n3. Feasibility study whether Appzen can flag it as medium or low risk if the users attach only credit card statement (Jan 7, 2022)`n"

# --- Row heights re-flow after the status/content edits --------------------
$ws.Rows.Item(3).RowHeight = 101.5
$ws.Rows.Item(10).RowHeight = 290
$ws.Rows.Item(13).RowHeight = 87
$ws.Rows.Item(14).RowHeight = 72.5
$ws.Rows.Item(27).RowHeight = 87
$ws.Rows.Item(28).RowHeight = 58
$ws.Rows.Item(29).RowHeight = 101.5
$ws.Rows.Item(30).RowHeight = 101.5
$ws.Rows.Item(32).RowHeight = 101.5
$ws.Rows.Item(33).RowHeight = 43.5

# --- Scroll position / selection -------------------------------------------
$ws.Activate() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 31
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.TopLeftCell = $ws.Range("C31")
} catch {
}
$ws.Range("F35").Select() | Out-Null
